$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. After the "ESTADO DEL ARTE" title paragraph, insert two new bold
#    paragraphs: "Intro " and a tab + "ANTECEDENTES".
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$insertPoint1 = $d.Range($p1.Range.End, $p1.Range.End)
$xmlFrag1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Intro</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ANTECEDENTES</w:t></w:r>
</w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertPoint1.InsertXML($xmlFrag1)

# The fragment above ends with a spare empty paragraph so that the
# "ANTECEDENTES" paragraph keeps its own bold paragraph mark instead of
# merging into the following "No hay proyectos..." paragraph. Remove
# that spare empty paragraph now (deleting it merges forward into
# "No hay proyectos...", which had no special paragraph-mark formatting
# to begin with, so nothing is lost).
$spacer = $d.Paragraphs.Item(4)
$spacer.Range.Delete()

# ---------------------------------------------------------------------
# 2. Fix the hyperlink text that was split across three runs
#    ("https://www.david-colso" + "n" + ".com/...") into one continuous
#    run of text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("https://www.david-colson.com/2020/02/09/making-a-simple-ecs.html", $true, $false, $false, $false, $false, $true, 1, $false, "https://www.david-colson.com/2020/02/09/making-a-simple-ecs.html", 2)

# ---------------------------------------------------------------------
# 3. At the end of the document, add the new "TECNOLOGÍAS PARA EL
#    DESARROLLO" section.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint2 = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$xmlFrag2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p>
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>TECNOLOGÍAS PARA EL DESARROLLO</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Las tecnologías que vamos a aplicar para la elaboración del proyecto son principalmente C++20, como la base de todo el motor ECS, complementado con un motor gráfico para obtener una demo sencilla en dos dimensiones. Este motor será </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Irrlicht</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Engine</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>.</w:t></w:r>
</w:p>
<w:p/>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertPoint2.InsertXML($xmlFrag2)

Write-Output "Edit complete."
